$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns L (Total penalty points?), M (Rank), N (Status/consequence)
# for each of the 10 student rows.

$data = @(
    @{ Row = 1;  L = 180; M = "Rank 1";  N = "No problem" },
    @{ Row = 2;  L = 180; M = "Rank 7";  N = "Warning" },
    @{ Row = 3;  L = 150; M = "Rank 5";  N = "Reprimand" },
    @{ Row = 4;  L = 120; M = "Rank 10"; N = "No problem" },
    @{ Row = 5;  L = 200; M = "Rank 3";  N = "Reprimand" },
    @{ Row = 6;  L = 100; M = "Rank 1";  N = "Warning" },
    @{ Row = 7;  L = 180; M = "Rank 5";  N = "No problem" },
    @{ Row = 8;  L = 160; M = "Rank 12"; N = "Warning" },
    @{ Row = 9;  L = 100; M = "Rank 12"; N = "No problem" },
    @{ Row = 10; L = 170; M = "Rank 17"; N = "Reprimand" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
}

# Set column N width to match new content column
$ws.Columns.Item(14).ColumnWidth = 12.29

# Update selection / view state to match the post-edit state
$ws.Range("I11").Select()
